# Update docx golden tests for style changes.
#
# 1. Add a new "Abstract Title" paragraph style (AbstractTitle), based on
#    Normal, followed by Abstract.
# 2. Abstract style: change space-before from 300 twips (15pt) to 100
#    twips (5pt); space-after stays 300 twips (15pt).
# 3. ImportTok character style: add green (008000) bold formatting.
# 4. BuiltInTok character style: add green (008000) formatting.

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style -----------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. Abstract style spacing tweak ----------------------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. ImportTok character style --------------------------------------
$importTok = $d.Styles.Item("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# --- 4. BuiltInTok character style --------------------------------------
$builtInTok = $d.Styles.Item("BuiltInTok")
$builtInTok.Font.Color = 32768

Write-Output "style updates applied"
